$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oct 15 sampling data: corrected reading in C10 (255 -> 253)
$ws.Range("C10").Value = 253

# Move the active selection to C10 (reflects where the user left off editing)
$ws.Range("C10").Select()

# Best-effort: reposition the document window to match the recorded view state.
# (No-op if the host doesn't expose window geometry through the object model.)
try {
    $win = $excel.ActiveWindow
    $win.Left = 720
    $win.Top = 4320
} catch {
}
